$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = 511367917
$ws.Range("G12").Select()
